$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "<alpha>"

$ws.Range("C3").Value = 23

$ws.Range("B4").Value = "<been>"
$ws.Range("C4").Value = 30

$ws.Range("C5").Value = 29

$ws.Range("B6").Value = "<or>"
$ws.Range("C6").Value = 27

$ws.Range("C7").Value = 31

$ws.Range("C8").Value = 28

$ws.Range("C10").Value = 27

$ws.Range("B11").Value = "<get>"
$ws.Range("C11").Value = 25

$ws.Range("B12").Value = "<in>"
$ws.Range("C12").Value = 22

$ws.Range("B13").Value = "<tango>"
$ws.Range("C13").Value = 30

$ws.Range("C14").Value = 22

$ws.Range("C15").Value = 24

$ws.Range("B16").Value = "<from>"
$ws.Range("C16").Value = 17
